$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 ("2021年") appended below the existing data (last row was 12, "2020年").
# Copy the formatting of the row-13 label cell from A12 (bold/bordered/centered style
# used throughout column A) so the new label cell matches the existing look.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"

$ws.Range("B13").Value = 59.4
$ws.Range("C13").Value = 62.9
$ws.Range("D13").Value = 9732.6
$ws.Range("F13").Value = 242150.5
$ws.Range("G13").Value = 242321
$ws.Range("H13").Value = 151717.4
$ws.Range("I13").Value = 81170.7
$ws.Range("J13").Value = 56698.9
$ws.Range("K13").Value = 208616.7
$ws.Range("L13").Value = 164924.8
$ws.Range("N13").Value = 526.9
$ws.Range("O13").Value = 49990.4
$ws.Range("P13").Value = 46691.4
$ws.Range("Q13").Value = 373541.5
$ws.Range("R13").Value = 1964.5
$ws.Range("S13").Value = 17678.3

# E13 and M13 stay blank (matching the source pattern of empty cells elsewhere
# in the same columns, e.g. E12/M12, K2/L2).

Write-Output "row 13 added"
